# Update December row (row 13) figures: Good, Reject, BU
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 2419
$ws.Range("C13").Value = 179
$ws.Range("D13").Value = 455

# Move/update the active selection shown when the sheet was last saved
$ws.Range("J20").Select()
